$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Mary): update password and login-status flags
$ws.Range("D2").Value = "password1"
$ws.Range("E2").Value = $true
$ws.Range("F2").Value = $false

# Row 3 (Caleb): update password and login-status flag
$ws.Range("D3").Value = "password!"
$ws.Range("F3").Value = $false
